$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2026-01-14 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-01-15 Thursday", 2)

# Update the division problems in the table by exact cell position
# (row, col) are 1-indexed as used by Word's Table.Cell(row, col) API.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "49÷7=" },
    @{ Row = 1;  Col = 2; Text = "50÷5=" },
    @{ Row = 1;  Col = 3; Text = "36÷3=" },
    @{ Row = 1;  Col = 4; Text = "30÷6=" },
    @{ Row = 1;  Col = 5; Text = "77÷6=" },

    @{ Row = 5;  Col = 1; Text = "88÷5=" },
    @{ Row = 5;  Col = 2; Text = "60÷9=" },
    @{ Row = 5;  Col = 3; Text = "32÷5=" },
    @{ Row = 5;  Col = 4; Text = "38÷6=" },
    @{ Row = 5;  Col = 5; Text = "26÷8=" },

    @{ Row = 9;  Col = 1; Text = "23÷9=" },
    @{ Row = 9;  Col = 2; Text = "41÷5=" },
    @{ Row = 9;  Col = 3; Text = "79÷4=" },
    @{ Row = 9;  Col = 4; Text = "40÷6=" },
    @{ Row = 9;  Col = 5; Text = "44÷4=" },

    @{ Row = 13; Col = 1; Text = "54÷8=" },
    @{ Row = 13; Col = 2; Text = "83÷4=" },
    @{ Row = 13; Col = 3; Text = "58÷7=" },
    @{ Row = 13; Col = 4; Text = "92÷4=" },
    @{ Row = 13; Col = 5; Text = "44÷5=" },

    @{ Row = 17; Col = 1; Text = "70÷5=" },
    @{ Row = 17; Col = 2; Text = "50÷3=" },
    @{ Row = 17; Col = 3; Text = "30÷4=" },
    @{ Row = 17; Col = 4; Text = "78÷2=" },
    @{ Row = 17; Col = 5; Text = "23÷9=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
